# Updated capital structure database
# Mexico Brokerage & Investment Banking sheet - refresh row 2-4 figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (aggregate row) ---
$ws.Range("D2").Value = -0.09250000000000003
$ws.Range("E2").Value = -0.11
$ws.Range("I2").Value = 0.0003773564202369911
$ws.Range("J2").Value = 0.0003217416639793559
$ws.Range("K2").Value = 36.24
$ws.Range("L2").Value = 0.1411764705882353
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("U2").Value = 62
$ws.Range("V2").Value = 0.04625139873181649
$ws.Range("W2").Value = 0.06056870800545333
$ws.Range("X2").Value = 0.09249044919190408
$ws.Range("Y2").Value = -0.03192174118645075
$ws.Range("Z2").Value = 0.0518510369443166
$ws.Range("AA2").Value = [double]"6.443901680652244e-05"
$ws.Range("AB2").Value = 0.03974737017055661
$ws.Range("AC2").Value = -0.0396829311537501
$ws.Range("AD2").Value = 1930
$ws.Range("AE2").Value = 4.020663034625822
$ws.Range("AF2").Value = 1934.020663034626
$ws.Range("AG2").Value = 1872.020663034626
$ws.Range("AH2").Value = 0.5906271060883438
$ws.Range("AI2").Value = 0.7784594178476573
$ws.Range("AJ2").Value = 0.5827264193427069
$ws.Range("AK2").Value = 0.7727892564660918
$ws.Range("AN2").Value = 2142.064372918979
$ws.Range("AP2").Value = 2077.714387385822

# --- Row 3: company is now Value Grupo Financiero ---
$ws.Range("B3").Value = "Value Grupo Financiero, S.A.B. de C.V. (BMV:VALUEGF O)"
$ws.Range("D3").Value = -0.291
$ws.Range("E3").Value = -0.3779999999999999
$ws.Range("I3").Value = 0.004967558619222339
$ws.Range("J3").Value = 0.004967558619222339
$ws.Range("K3").Value = 4.24
$ws.Range("L3").Value = 0.2174358974358974
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("U3").Value = 8.199999999999999
$ws.Range("V3").Value = 0.00765068109721963
$ws.Range("W3").Value = 0.01645962732919255
$ws.Range("X3").Value = 0.04131730443043367
$ws.Range("Y3").Value = -0.02485767710124113
$ws.Range("Z3").Value = 0.02594393815793974
$ws.Range("AA3").Value = 0.0001288780336130449
$ws.Range("AB3").Value = 0.03877452372730188
$ws.Range("AC3").Value = -0.03864564569368884
$ws.Range("AD3").Value = 340.5
$ws.Range("AE3").Value = 4.020663034625822
$ws.Range("AF3").Value = 344.5206630346258
$ws.Range("AG3").Value = 336.3206630346259
$ws.Range("AH3").Value = 0.2432504672327887
$ws.Range("AI3").Value = 0.5896088994104427
$ws.Range("AJ3").Value = 0.2388436388042377
$ws.Range("AK3").Value = 0.583767749733379
$ws.Range("AN3").Value = 377.9134295227525
$ws.Range("AP3").Value = 373.2748757321041

# --- Row 4: company is now Corporacion Actinver ---
$ws.Range("B4").Value = "Corporación Actinver, S. A. B. de C. V. (BMV:ACTINVR B)"
$ws.Range("D4").Value = 0.106
$ws.Range("E4").Value = 0.158
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 32
$ws.Range("L4").Value = 0.1349072512647555
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("U4").Value = 53.8
$ws.Range("V4").Value = 0.2002232973576479
$ws.Range("W4").Value = 0.1046777886817141
$ws.Range("X4").Value = 0.1436635939533745
$ws.Range("Y4").Value = -0.03898580527166036
$ws.Range("Z4").Value = 0.05648829511085708
$ws.Range("AA4").Value = 0
$ws.Range("AB4").Value = 0.04072021661381135
$ws.Range("AC4").Value = -0.04072021661381135
$ws.Range("AD4").Value = 1589.5
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 1589.5
$ws.Range("AG4").Value = 1535.7
$ws.Range("AH4").Value = 0.8553976966957271
$ws.Range("AI4").Value = 0.8365349192147782
$ws.Range("AJ4").Value = 0.8510862336510752
$ws.Range("AK4").Value = 0.8317716514109299

# --- Columns/cells no longer present in the refreshed dataset ---
$ws.Range("T2").ClearContents()
$ws.Range("T3").ClearContents()
$ws.Range("T4").ClearContents()
$ws.Range("AN4").ClearContents()
$ws.Range("AP4").ClearContents()
